$d = $word.ActiveDocument

$replacements = @(
    @{old="34×83="; new="54×55="},
    @{old="86×99="; new="28×98="},
    @{old="68×73="; new="34×26="},
    @{old="78×25="; new="61×35="},
    @{old="87×52="; new="88×75="},
    @{old="78×70="; new="70×43="},
    @{old="37×37="; new="53×97="},
    @{old="35×16="; new="95×53="},
    @{old="59×29="; new="52×34="},
    @{old="49×89="; new="81×76="},
    @{old="47×84="; new="63×48="},
    @{old="48×98="; new="97×88="},
    @{old="14×21="; new="60×53="},
    @{old="22×16="; new="79×76="},
    @{old="30×25="; new="19×36="},
    @{old="38×77="; new="12×62="},
    @{old="47×69="; new="72×17="},
    @{old="75×18="; new="11×85="},
    @{old="43×45="; new="23×68="},
    @{old="54×67="; new="49×79="},
    @{old="97×27="; new="95×74="},
    @{old="78×61="; new="57×94="},
    @{old="12×17="; new="38×58="},
    @{old="69×69="; new="92×94="},
    @{old="90×36="; new="82×25="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
